# Auto update Excel log 2026-02-04 14:29:13
# Appends the latest batch of sensor readings to the PIR, Humidity and
# Temperature logs (all three sheets share the same Date/Timestamp/Hour/
# Location/Value/Status layout).

$wb = $excel.ActiveWorkbook

# --- PIR sheet: rows 324-337 -------------------------------------------------
$pirData = @(
    @("2026-02-04","14:28:09","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:28:10","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:28:14","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:28:15","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:28:23","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:28:27","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:28:28","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:28:35","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:28:41","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:28:46","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:28:50","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:28:53","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:29:01","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:29:06","14:00","Bathroom","No Motion","Inactive")
)

$wsPIR = $wb.Worksheets.Item("PIR")
$r = 324
foreach ($row in $pirData) {
    # Column A holds a "YYYY-MM-DD" string; format as Text first so Excel
    # keeps it as literal text instead of auto-converting it to a date serial.
    $cellA = $wsPIR.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]

    $wsPIR.Cells.Item($r, 2).Value = $row[1]
    $wsPIR.Cells.Item($r, 3).Value = $row[2]
    $wsPIR.Cells.Item($r, 4).Value = $row[3]
    $wsPIR.Cells.Item($r, 5).Value = $row[4]
    $wsPIR.Cells.Item($r, 6).Value = $row[5]

    $r++
}

# --- Humidity sheet: rows 266-277 -------------------------------------------
$humData = @(
    @("2026-02-04","14:28:07","14:00","Bathroom","79.4%","Active"),
    @("2026-02-04","14:28:11","14:00","Bathroom","79.3%","Active"),
    @("2026-02-04","14:28:16","14:00","Bathroom","78.5%","Active"),
    @("2026-02-04","14:28:21","14:00","Bathroom","79.4%","Active"),
    @("2026-02-04","14:28:31","14:00","Bathroom","79.5%","Active"),
    @("2026-02-04","14:28:36","14:00","Bathroom","78.6%","Active"),
    @("2026-02-04","14:28:42","14:00","Bathroom","79.5%","Active"),
    @("2026-02-04","14:28:47","14:00","Bathroom","78.5%","Active"),
    @("2026-02-04","14:28:51","14:00","Bathroom","79.5%","Active"),
    @("2026-02-04","14:28:56","14:00","Bathroom","78.6%","Active"),
    @("2026-02-04","14:29:02","14:00","Bathroom","79.6%","Active"),
    @("2026-02-04","14:29:07","14:00","Bathroom","78.7%","Active")
)

$wsHumidity = $wb.Worksheets.Item("Humidity")
$r = 266
foreach ($row in $humData) {
    $cellA = $wsHumidity.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]

    $wsHumidity.Cells.Item($r, 2).Value = $row[1]
    $wsHumidity.Cells.Item($r, 3).Value = $row[2]
    $wsHumidity.Cells.Item($r, 4).Value = $row[3]

    # Column E holds a percentage-looking string (e.g. "79.4%") which Excel
    # would otherwise interpret as a numeric percent value - keep as text.
    $cellE = $wsHumidity.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $row[4]

    $wsHumidity.Cells.Item($r, 6).Value = $row[5]

    $r++
}

# --- Temperature sheet: rows 266-276 ----------------------------------------
$tempData = @(
    @("2026-02-04","14:28:08","14:00","Bathroom","24.4C","Active"),
    @("2026-02-04","14:28:12","14:00","Bathroom","24.3C","Active"),
    @("2026-02-04","14:28:17","14:00","Bathroom","24.3C","Active"),
    @("2026-02-04","14:28:22","14:00","Bathroom","24.3C","Active"),
    @("2026-02-04","14:28:32","14:00","Bathroom","24.3C","Active"),
    @("2026-02-04","14:28:37","14:00","Bathroom","24.4C","Active"),
    @("2026-02-04","14:28:42","14:00","Bathroom","24.3C","Active"),
    @("2026-02-04","14:28:47","14:00","Bathroom","24.3C","Active"),
    @("2026-02-04","14:28:52","14:00","Bathroom","24.4C","Active"),
    @("2026-02-04","14:28:57","14:00","Bathroom","24.3C","Active"),
    @("2026-02-04","14:29:02","14:00","Bathroom","24.3C","Active")
)

$wsTemperature = $wb.Worksheets.Item("Temperature")
$r = 266
foreach ($row in $tempData) {
    $cellA = $wsTemperature.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]

    $wsTemperature.Cells.Item($r, 2).Value = $row[1]
    $wsTemperature.Cells.Item($r, 3).Value = $row[2]
    $wsTemperature.Cells.Item($r, 4).Value = $row[3]
    $wsTemperature.Cells.Item($r, 5).Value = $row[4]
    $wsTemperature.Cells.Item($r, 6).Value = $row[5]

    $r++
}
